# Removed Test Case Inter-Dependency
# Replace the hard-coded product/short names on the ProductLoanInput sheet
# with values that don't collide with other automated test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B1").Value = "2609-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-ADD-VAR-INST-OVERDUE-FEE-1st"
$ws.Range("B2").Value = "260a"

$ws.Range("B3").Select() | Out-Null
